# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 44

# New header cells, formatted like the existing header row (bold, bordered,
# centered) by copying the format from the neighboring header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# The 2005 Mets finished 83-79-0; stamp that record onto every player row.
$ws.Range("AD2:AD$lastRow").Value = 83
$ws.Range("AE2:AE$lastRow").Value = 79
$ws.Range("AF2:AF$lastRow").Value = 0
